$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.001.12"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.183.15"
$ws.Range("E3").Value = "  -4.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "589.93"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.42%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "134.85"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.02%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.181.33"
$ws.Range("E8").Value = "  -4.05%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.516"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  -6.13%  "
$ws.Range("E11").Value = "  -5.86%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.453"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("E13").Value = "  -5.00%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.62"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "3.707.63"
$ws.Range("E15").Value = "  -4.04%  "
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "3.181.39"
$ws.Range("E17").Value = "  -4.08%  "
$ws.Range("D18").Value = "63.003.10"
$ws.Range("E18").Value = "  -1.37%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.56"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.43%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "461.53"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.25%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.03"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E22").Value = "  -6.27%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.63"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -4.57%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "13.35"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -4.83%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "82.64"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.03%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.77"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.82%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.69"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -6.72%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.84%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "27.14"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -6.30%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.103"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.53%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.21%  "
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("E36").Value = "  -4.63%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "51.32"
$c.Style = "Normal"
$ws.Range("D38").Value = "0.0₃0706"
$ws.Range("E38").Value = "  -5.59%  "
$ws.Range("E39").Value = "  -3.08%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "405.15"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.90%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "8.07"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.45%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.64"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.07%  "
$ws.Range("E43").Value = "  -5.72%  "
$ws.Range("D44").Value = "2.795.37"
$ws.Range("E44").Value = "  -10.45%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.252"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.92%  "
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "34.70"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -5.99%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "25.17"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("E51").Value = "  -1.91%  "
